$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at O and P (old O shifts to Q, etc.)
$ws.Columns.Item(15).Insert()
$ws.Columns.Item(15).Insert()

# Update header text for M1 and N1 (renamed, no column shift)
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Set header text for the two newly inserted columns O1 and P1
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Fill the new O and P columns (rows 2 through 51) with 0
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = 0
    $ws.Cells.Item($r, 16).Value = 0
}

Write-Host "Done"
